$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill User Story (C) and Sprint (D) columns for rows 2-7
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "sprint 2"

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "sprint 2"

$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "sprint 2"

$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "sprint 2"

$ws.Range("C6").Value = 19
$ws.Range("D6").Value = "sprint 2"

$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "sprint 2"

# Adjust row 1 height
$ws.Rows.Item(1).RowHeight = 18

# Update the current selection to E7
$ws.Range("E7").Select()
